# Commitment import now updates fields
# Inserts a new "Onboarding Completed" column (I) before the existing
# Fund Close column, shifting Fund Close..As Of from I..N to J..O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I, shifting existing I:N -> J:O
$ws.Range("I1:I7").EntireColumn.Insert()

# Header for new column
$ws.Range("I1").Value = "Onboarding Completed"

# Values for rows 2-7 (row 3 = "No", rest = "Yes")
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "No"
$ws.Range("I4").Value = "Yes"
$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"

# Clear the stray style on the shifted header cells (previously s="4",
# a duplicate of the default style) so they fall back to the default
$ws.Range("L1:O1").Style = "Normal"

# Selection matches the author's final click target
$ws.Range("I1:I7").Select()
